$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Capture current header values for columns C..G (3..7)
    $c1 = $ws.Cells.Item(1, 3).Value()
    $d1 = $ws.Cells.Item(1, 4).Value()
    $e1 = $ws.Cells.Item(1, 5).Value()
    $f1 = $ws.Cells.Item(1, 6).Value()
    $g1 = $ws.Cells.Item(1, 7).Value()

    # Rotate left: C<-D, D<-E, E<-F, F<-G, G<-C (old)
    $ws.Cells.Item(1, 3).Value = $d1
    $ws.Cells.Item(1, 4).Value = $e1
    $ws.Cells.Item(1, 5).Value = $f1
    $ws.Cells.Item(1, 6).Value = $g1
    $ws.Cells.Item(1, 7).Value = $c1

    # Column widths follow the header text that now lives in each column
    # (Excel recomputes "best fit" widths after the header swap).
    $ws.Columns.Item(3).ColumnWidth = 7.166666666666667
    $ws.Columns.Item(7).ColumnWidth = 8.666666666666666
}
